$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.405.74'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '''2.017.64'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +5.95%  '
$ws.Range('D5').Value = '''325.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''0.5138'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.07%  '
$ws.Range('D8').Value = '''0.4211'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.28%  '
$ws.Range('D9').Value = '''0.08730'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.90%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '''43.58'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.89%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '''1.135'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.78%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value = '''24.83'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.62%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''2.012.58'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.38%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''6.599'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.01%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '''7.479'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.16%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').Value = '''1.004'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '''94.64'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.12%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.00001114'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.09%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '''0.06531'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '''18.98'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.54%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '''1.001'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '''6.207'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.65%  '
$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').Value = '''30.461.60'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.75%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '''11.86'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.76%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''2.240'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.32%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '''2.251.54'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +5.99%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''22.41'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''162.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '''2.427'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.51%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''131.53'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.140'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.76%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.1054'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.08%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''6.070'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.830'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''1.374'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +15.60%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').Value = '''0.02531'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.26%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '''5.482'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.26%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.06670'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.41%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '''12.33'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +9.25%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = '''0.2201'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''9.074'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +5.27%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '''0.6662'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.15%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''1.233'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.60%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '''1.000'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''13.66'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.36%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.6181'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.29%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '''2.191'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '''3.663'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '''1.265'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.14%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''124.58'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''81.01'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.47%  '
